$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Mark the final "footer" paragraph as English (US) text. Do this before
#    touching the table below (indexing into Paragraphs stays document-order
#    stable that way).
# ---------------------------------------------------------------------------
$last = $d.Content.Paragraphs.Item($d.Content.Paragraphs.Count)
$last.Range.LanguageID = "en-US"

# ---------------------------------------------------------------------------
# 2) Table formatting: give the "Key Highlights" table an explicit auto width
#    and a standard tblLook (heading row + first column banding, no column
#    banding) to match what Word writes when a table style is (re)applied.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

$t.PreferredWidthType = 1   # wdPreferredWidthAuto -> <w:tblW w:type="auto" .../>
$t.PreferredWidth = 0

$t.ApplyStyleHeadingRows  = $true
$t.ApplyStyleLastRow      = $false
$t.ApplyStyleFirstColumn  = $true
$t.ApplyStyleLastColumn   = $false
$t.ApplyStyleRowBands     = $true
$t.ApplyStyleColumnBands  = $false

# ---------------------------------------------------------------------------
# 3) Give the document an explicit section (page size / margins / columns)
#    matching a standard A4 layout instead of relying on implicit defaults.
# ---------------------------------------------------------------------------
$section = $d.Sections.Item(1)
$ps = $section.PageSetup

$ps.PageWidth  = 595.3
$ps.PageHeight = 841.9

$ps.TopMargin    = 70.85
$ps.BottomMargin = 56.7
$ps.LeftMargin   = 70.85
$ps.RightMargin  = 70.85
$ps.Gutter       = 0

$ps.HeaderDistance = 35.4
$ps.FooterDistance = 35.4

$ps.TextColumns.Spacing = 35.4
